# Adding new test cases for watch list (TestCase_E30, TestCase_E31, TestCase_E32)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 28 (TestCase_E27) gains a third linked Jira issue / requirement ---
$ws.Range("B28").Value = "OPQA-314 ||OPQA-317" + [char]10 + "||OPQA-327"
$ws.Range("C28").Value = "Verify that user is able to name the watchlists||Verify that a user can add description to his watchlist||Verify that watchlist name is customizable"
$ws.Rows.Item(28).RowHeight = 45

# --- Append three brand-new test cases as rows 31-33 ---
# Start from row 29's cell formatting (bordered, description column wrapped) as a template.
$ws.Range("A29:E29").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$ws.Range("A32:E32").PasteSpecial(-4122)
$ws.Range("A33:E33").PasteSpecial(-4122)

# Column B on these new rows is a plain bordered cell (no wrap), like column A/D/E.
$ws.Range("B2").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B33").PasteSpecial(-4122)

# Row 32's Description cell keeps the shaded/filled wrap style (matches row 30's C column).
$ws.Range("C30").Copy()
$ws.Range("C32").PasteSpecial(-4122)

$ws.Range("A31").Value = "TestCase_E30"
$ws.Range("B31").Value = "OPQA-324"
$ws.Range("C31").Value = "Verify that a user has 1 watchlist by default once we try to watch an item"
$ws.Range("D31").Value = "Y"
$ws.Range("E31").Value = "PASS"

$ws.Range("A32").Value = "TestCase_E31"
$ws.Range("B32").Value = "OPQA-326"
$ws.Range("C32").Value = "Verify that user is able to have a watchlist with 0 item under it"
$ws.Range("D32").Value = "Y"
$ws.Range("E32").Value = "PASS"

$ws.Range("A33").Value = "TestCase_E32"
$ws.Range("B33").Value = "OPQA-328"
$ws.Range("C33").Value = "Verify that every user watchlist is private by default"
$ws.Range("D33").Value = "Y"
$ws.Range("E33").Value = "PASS"

# --- Match the final selection / active cell shown in the workbook ---
$ws.Range("E33").Select()
